$d = $word.ActiveDocument

# 1. Insert a new paragraph at the very beginning containing the title "test3"
#    centered, bold, size 14pt (sz 28 half-points).
$firstPara = $d.Paragraphs.First
$firstPara.Range.InsertParagraphBefore()

$titlePara = $d.Paragraphs.First
$titlePara.Alignment = 1

$titleRange = $titlePara.Range
$titleRange.Text = "test3"

$titleTextRange = $d.Range($titlePara.Range.Start, $titlePara.Range.Start + 5)
$titleTextRange.Font.Bold = 1
$titleTextRange.Font.Size = 14

# 2. Apply text corrections inside the big paragraphs.

$rng = $d.Content
$rng.Find.Execute("del mundo y viene alguien", $true, $false, $false, $false, $false, $true, 1, $false, "del mundo, y viene alguien", 2)

$rng = $d.Content
$rng.Find.Execute("atender de otros y dice la Biblia va hasta el de. a cada día su propio pan. y si nos ponemos", $true, $false, $false, $false, $false, $true, 1, $false, "atender de otros? y dice la Biblia: basta el de a cada día su propio afán y si nos ponemos", 2)

$rng = $d.Content
$rng.Find.Execute("nuestras preocupaciones no son acerca", $true, $false, $false, $false, $false, $true, 1, $false, "nuestras preocupaciones-, no son acerca", 2)

$rng = $d.Content
$rng.Find.Execute("capítulo 3, verso 25: lo que temí me aconteció. para el que cree, en realidad, todo es posible. si tú piensas: no, yo no voy a poder, no vas a poder.", $true, $false, $false, $false, $false, $true, 1, $false, "capítulo 3, verso 25, job decía: lo que temí aconteció para el que cree en realidad, todo es posible. si tú piensas no, yo no voy a poder, no vas a poder.", 2)

$rng = $d.Content
$rng.Find.Execute("nos seguimos adelante. sino que estamos concentrados en esa sola cosa, en esa sola cosa, en esa sola cosa le damos vueltas", $true, $false, $false, $false, $false, $true, 1, $false, "nos seguimos adelante sino que estamos concentrados en esa sola cosa. en esa sola cosa, en esa sola cosa, le damos vueltas", 2)

Write-Host "Done"
